# This NATMI LR-pairs sheet recomputes ligand/receptor expression (TPM)
# statistics per "Sending cluster" (column A) / "Target cluster" (column D)
# pair and derives the edge-weight / specificity columns from them.
#
# New per-cluster ligand values (average/total expression), keyed by the
# Sending cluster name (column A):
$newLigand = @{
    "ECs"            = @(0.225007, 0.675021)
    "FAPs"           = @(4.404016666666667, 13.21205)
    "MuSCs"          = @(0.2943246666666666, 0.8829739999999999)
    "Resolving-Mac"  = @(1.060660666666667, 3.181982)
}

# New per-cluster receptor values (# expressing cells, detection rate,
# average expression, total expression), keyed by the Target cluster name
# (column D):
$newReceptor = @{
    "ECs"            = @(3, 1,                    7.107333666666666,    21.322001)
    "FAPs"           = @(3, 1,                    1.627877666666667,    4.883633)
    "MuSCs"          = @(3, 1,                    0.8135026666666666,   2.440508)
    "Resolving-Mac"  = @(2, 0.6666666666666666,   0.09009266666666667,  0.270278)
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row   # xlUp
if ($lastRow -lt $firstRow) { $lastRow = 17 }

# --- Pass 1: update the directly-measured columns (G, H, K, L, M, N) ---
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $sending = $ws.Cells.Item($r, 1).Value2
    $target  = $ws.Cells.Item($r, 4).Value2

    $lig = $newLigand[$sending]
    $ws.Cells.Item($r, 7).Value2  = $lig[0]   # G: Ligand average expression value
    $ws.Cells.Item($r, 8).Value2  = $lig[1]   # H: Ligand total expression value

    $rec = $newReceptor[$target]
    $ws.Cells.Item($r, 11).Value2 = $rec[0]   # K: Receptor-expressing cells
    $ws.Cells.Item($r, 12).Value2 = $rec[1]   # L: Receptor detection rate
    $ws.Cells.Item($r, 13).Value2 = $rec[2]   # M: Receptor average expression value
    $ws.Cells.Item($r, 14).Value2 = $rec[3]   # N: Receptor total expression value
}

# --- Pass 2: derive specificity (I, J, O, P) from sums over unique clusters ---
$sumG = 0.0
$sumH = 0.0
foreach ($v in $newLigand.Values) { $sumG += $v[0]; $sumH += $v[1] }

$sumM = 0.0
$sumN = 0.0
foreach ($v in $newReceptor.Values) { $sumM += $v[2]; $sumN += $v[3] }

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $g = $ws.Cells.Item($r, 7).Value2
    $h = $ws.Cells.Item($r, 8).Value2
    $m = $ws.Cells.Item($r, 13).Value2
    $n = $ws.Cells.Item($r, 14).Value2

    $ws.Cells.Item($r, 9).Value2  = $g / $sumG   # I: Ligand derived specificity of average expression value
    $ws.Cells.Item($r, 10).Value2 = $h / $sumH   # J: Ligand derived specificity of total expression value
    $ws.Cells.Item($r, 15).Value2 = $m / $sumM   # O: Receptor derived specificity of average expression value
    $ws.Cells.Item($r, 16).Value2 = $n / $sumN   # P: Receptor derived specificity of total expression value

    $ws.Cells.Item($r, 17).Value2 = $g * $m      # Q: Edge average expression weight
    $ws.Cells.Item($r, 18).Value2 = $h * $n      # R: Edge total expression weight
}

# --- Pass 3: derive edge specificity (S, T) from sums over all edges ---
$sumQ = 0.0
$sumR = 0.0
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $sumQ += $ws.Cells.Item($r, 17).Value2
    $sumR += $ws.Cells.Item($r, 18).Value2
}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $q = $ws.Cells.Item($r, 17).Value2
    $rr = $ws.Cells.Item($r, 18).Value2
    $ws.Cells.Item($r, 19).Value2 = $q / $sumQ   # S: Edge average expression derived specificity
    $ws.Cells.Item($r, 20).Value2 = $rr / $sumR  # T: Edge total expression derived specificity
}
